# forecast_series_EXPORT_yoy_AR2_50_9_full.xlsx - bugfixed evaluation / simulated rt_data
#
# A new (oldest) observation row is inserted at row 2, pushing all existing rows
# down by one (so the sheet grows from 52 to 53 data rows). The recomputed AR(2)
# forecast bugfix also changes most of the y_0_forecast (C) / y_1_forecast (E)
# values throughout the series, so every row is rewritten explicitly below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 2:52 down to 3:53 to make room for the new first data row.
$ws.Rows("2:2").Insert()

# Row insertion copies row 1's (header) formatting onto the new row; restore the
# plain per-row look used by the rest of the sheet (date format only on column A).
$ws.Range("B2:E2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 2
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 11.13090654781821
$ws.Range("D2").Value = 2008
$ws.Range("E2").Clear()

# row 3
$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 2008
$ws.Range("C3").Clear()
$ws.Range("D3").Value = 2009
$ws.Range("E3").Clear()

# row 4
$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = 2008
$ws.Range("C4").Value = 4.672550446571067
$ws.Range("D4").Value = 2009
$ws.Range("E4").Clear()

# row 5
$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = 2009
$ws.Range("C5").Clear()
$ws.Range("D5").Value = 2010
$ws.Range("E5").Clear()

# row 6
$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = 2009
$ws.Range("C6").Value = -14.45332333832743
$ws.Range("D6").Value = 2010
$ws.Range("E6").Clear()

# row 7
$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = 2010
$ws.Range("C7").Value = 4.317520552235576
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = 3.856664261949816

# row 8
$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 2010
$ws.Range("C8").Value = 8.600536527919633
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = 9.002271992040312

# row 9
$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 2011
$ws.Range("C9").Value = 10.48563750975209
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 5.658864198748459

# row 10
$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2011
$ws.Range("C10").Value = 10.25770250047622
$ws.Range("D10").Value = 2012
$ws.Range("E10").Value = 7.550992341868956

# row 11
$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 2012
$ws.Range("C11").Value = 4.65880603412161
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = 5.409232631930561

# row 12
$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 2012
$ws.Range("C12").Value = 4.639893381363169
$ws.Range("D12").Value = 2013
$ws.Range("E12").Value = 5.799303245920884

# row 13
$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = 2013
$ws.Range("C13").Value = -0.8693696108860949
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = 2.311254283099773

# row 14
$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2013
$ws.Range("C14").Value = 0.3058963467304165
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = 2.638010271840896

# row 15
$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 2014
$ws.Range("C15").Value = 4.057580120451165
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = 3.242274909585041

# row 16
$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 2014
$ws.Range("C16").Value = 4.068173739091874
$ws.Range("D16").Value = 2015
$ws.Range("E16").Value = 5.156937396706884

# row 17
$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 2015
$ws.Range("C17").Value = 3.942841799197594
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = 3.826653192455631

# row 18
$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 4.984288257750213
$ws.Range("D18").Value = 2016
$ws.Range("E18").Value = 3.4064284328156

# row 19
$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = 2016
$ws.Range("C19").Value = 2.435545128806416
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = 4.047617271894799

# row 20
$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 2016
$ws.Range("C20").Value = 1.878184267712912
$ws.Range("D20").Value = 2017
$ws.Range("E20").Value = 2.129835064860464

# row 21
$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 2017
$ws.Range("C21").Value = 4.218672000695523
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = 4.02808863798465

# row 22
$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2017
$ws.Range("C22").Value = 4.695933104194339
$ws.Range("D22").Value = 2018
$ws.Range("E22").Value = 5.022591279638045

# row 23
$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 7.356933926419673
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 4.456100749631386

# row 24
$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = 4.174017263680696
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = 2.074691389445271

# row 25
$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = 2018
$ws.Range("C25").Value = 5.32644934790627
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = 3.123778338720062

# row 26
$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 2018
$ws.Range("C26").Value = 4.892602738886098
$ws.Range("D26").Value = 2019
$ws.Range("E26").Value = 0.6513682883433347

# row 27
$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = 1.443764323860086
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = 2.608754752952369

# row 28
$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = 1.699348375745302
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = 3.362240252406901

# row 29
$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 2019
$ws.Range("C29").Value = 0.4648465628725118
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = 0.4030385202844711

# row 30
$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = 2019
$ws.Range("C30").Value = 0.8049382522247184
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = 3.036929265763488

# row 31
$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = 0.642635895824295
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = 2.24075411271687

# row 32
$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = -2.856524424985296
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = 0.1239622353166103

# row 33
$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = 2020
$ws.Range("C33").Value = -9.810777096850787
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = -4.291312636561706

# row 34
$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 2020
$ws.Range("C34").Value = -8.784173899737169
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = 1.573231731123359

# row 35
$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = 5.992199201897175
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = 3.59700023027214

# row 36
$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = 6.240787792289715
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = 4.112367048326182

# row 37
$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 2021
$ws.Range("C37").Value = 5.797134106720514
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = 4.062819790119954

# row 38
$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = 5.110501195359984
$ws.Range("D38").Value = 2022
$ws.Range("E38").Value = 1.88131082127776

# row 39
$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = 5.158753599182209
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = 3.207070318919869

# row 40
$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = 3.815916106066686
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = 2.762663830671319

# row 41
$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 2022
$ws.Range("C41").Value = 4.232564748995715
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = 2.257990710759383

# row 42
$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2022
$ws.Range("C42").Value = 5.120680133083599
$ws.Range("D42").Value = 2023
$ws.Range("E42").Value = 2.522545412785848

# row 43
$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = 1.135459337362521
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = 2.746436727258894

# row 44
$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = 0.757583445265464
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = 2.564429185896056

# row 45
$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = 2023
$ws.Range("C45").Value = 0.08070151925247959
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = 1.796717073915977

# row 46
$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = -0.5532735011319234
$ws.Range("D46").Value = 2024
$ws.Range("E46").Value = -0.657715646732393

# row 47
$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = -2.103425609777143
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = 2.179038063425076

# row 48
$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = -0.152046383567539
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 2.659925441240518

# row 49
$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = 2024
$ws.Range("C49").Value = -0.9685570952743805
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = 1.778332450996523

# row 50
$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 2024
$ws.Range("C50").Value = -1.069674659641462
$ws.Range("D50").Value = 2025
$ws.Range("E50").Value = 0.5636794832278413

# row 51
$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = -2.92939668194816
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = 1.921265997709742

# row 52
$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = -1.93082584212636
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = 2.2373336846083

# row 53
$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = -2.436529450546909
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = 1.970944416887122
